# Update odds values on the active worksheet (sheet1) as described by the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 1.95
$ws.Range("I2").Value = 4.5
$ws.Range("S2").Value = 2.5
$ws.Range("T2").Value = 1.5

# Row 4
$ws.Range("G4").Value = 1.67
$ws.Range("H4").Value = 3.4
$ws.Range("I4").Value = 5.75
$ws.Range("J4").Value = 2.38
$ws.Range("K4").Value = 2.05
$ws.Range("M4").Value = 1.1
$ws.Range("N4").Value = 7
$ws.Range("Q4").Value = 1.78
$ws.Range("R4").Value = 2.1
$ws.Range("W4").Value = 4.33
$ws.Range("X4").Value = 1.2
$ws.Range("Y4").Value = 1.5
$ws.Range("Z4").Value = 2.5

# Row 5
$ws.Range("N5").Value = 8
$ws.Range("Q5").Value = 1.69
$ws.Range("R5").Value = 2.13
$ws.Range("U5").Value = 3.4
$ws.Range("V5").Value = 1.3
$ws.Range("AB5").Value = 1.67

# Row 6
$ws.Range("G6").Value = 3.6
$ws.Range("H6").Value = 3.1
$ws.Range("I6").Value = 2.15
$ws.Range("J6").Value = 4.33
$ws.Range("L6").Value = 3
$ws.Range("AA6").Value = 2.1
$ws.Range("AB6").Value = 1.67

# Row 7
$ws.Range("G7").Value = 3.05
$ws.Range("H7").Value = 2.87
$ws.Range("I7").Value = 2.42
$ws.Range("K7").Value = 1.95
$ws.Range("O7").Value = 1.4
$ws.Range("P7").Value = 2.5
$ws.Range("S7").Value = 2.18
$ws.Range("X7").Value = 1.2
$ws.Range("Y7").Value = 1.45
$ws.Range("Z7").Value = 2.37
$ws.Range("AA7").Value = 1.83
$ws.Range("AB7").Value = 1.78
